# Update "想去人数" (number of people interested) in column F for the
# sheets "展览" (index 1) and "全部类型" (index 4), which hold identical data.
# Values scraped/refreshed as of commit 456a3b4.

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 743
    4  = 259
    8  = 1679
    9  = 6160
    10 = 478
    12 = 292
    13 = 92
    14 = 369
    15 = 135
    16 = 6299
    18 = 1276
    20 = 115
    24 = 101
    26 = 9
    27 = 95
    30 = 82
    33 = 44
}

$sheetIndexes = @(1, 4)

foreach ($sheetIndex in $sheetIndexes) {
    $ws = $wb.Worksheets.Item($sheetIndex)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
